# Fill in the missing "75+" baseline count and move the active selection
# down to the next row, matching the author's manual data-entry workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 12

$ws.Range("B8").Select()
